$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine last used row (data goes from row 2 to row 258, header in row 1)
$lastRow = 258

# Add header for new column H
$ws.Cells.Item(1, 8).Value = "id mieszkania"

# Fill column H with sequential "id mieszkania" values: row 2 -> 1, row 3 -> 2, ... row 258 -> 257
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = $r - 1
}

# Update the view: scroll to bottom of data and select the new id column
$ws.Application.Goto($ws.Range("A246"), $true)
$ws.Range("H2:H258").Select()
